# Apply updates to the "Dados" sheet: correct several Idade (E column)
# values (decrement by 1) and swap the Filhos (B column) values for
# rows 27 and 28.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dados")

# Column E (Idade) corrections - decrement by 1
$eRows = @(5, 6, 11, 12, 13, 17, 18, 19, 20, 22, 28, 29, 30, 31, 34, 35, 36)
foreach ($r in $eRows) {
    $current = $ws.Cells.Item($r, 5).Value2
    $ws.Cells.Item($r, 5).Value2 = $current - 1
}

# Swap B27 and B28 (Filhos column) values
$ws.Cells.Item(27, 2).Value2 = 1
$ws.Cells.Item(28, 2).Value2 = 0
